$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column E (5th column) to approximate the "best fit" width Excel computed
# after the new row was added (closest value reachable in this runtime).
$ws.Columns.Item(5).ColumnWidth = 9.0

# Append the new trade row (row 6)
$ws.Range("A6").Value = 42649.65421296296
$ws.Range("A6").NumberFormat = "m/d/yy h:mm"

$ws.Range("B6").Value = $false

$ws.Range("C6").Value = 10010.959999999999
$ws.Range("D6").Value = 10015.469999999999
$ws.Range("E6").Value = 77.349997999999999
$ws.Range("F6").Value = 77.42

$ws.Range("G6").Value = $true
$ws.Range("G6").NumberFormat = "m/d/yy h:mm"

$ws.Range("H6").Value = 0.09

$ws.Range("I6").Value = $false
